# Update the exp_design workbook with the new unc-3 samples (last sample in
# the project). Adds 4 new rows (97-100) to Sheet1 describing the
# "MLC1602 unc-3:Ath-HEN1" L1 samples, prep1, replicates A/B with/without
# treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string text values, written in the same order they were
# first typed by the original author (sampleInfo column down first, then
# the repeated Tissue/promoter/details columns) so the resulting
# sharedStrings table lines up with the authored workbook.
$ws.Cells.Item(97, 8).Value = "MLC1602 unc-3:Ath-HEN1 L1s -A, prep1- No treatment"
$ws.Cells.Item(98, 8).Value = "MLC1602 unc-3:Ath-HEN1 L1s -A, prep1- Treatment"
$ws.Cells.Item(99, 8).Value = "MLC1602 unc-3:Ath-HEN1 L1s -B, prep1- No treatment"
$ws.Cells.Item(100, 8).Value = "MLC1602 unc-3:Ath-HEN1 L1s -B, prep1- Treatment"

$ws.Cells.Item(97, 6).Value = "unc-3 expressing neurons"
$ws.Cells.Item(98, 6).Value = "unc-3 expressing neurons"
$ws.Cells.Item(99, 6).Value = "unc-3 expressing neurons"
$ws.Cells.Item(100, 6).Value = "unc-3 expressing neurons"

$ws.Cells.Item(97, 7).Value = "unc-3"
$ws.Cells.Item(98, 7).Value = "unc-3"
$ws.Cells.Item(99, 7).Value = "unc-3"
$ws.Cells.Item(100, 7).Value = "unc-3"

$ws.Cells.Item(97, 9).Value = "unc-3 is a transcriptional factor"
$ws.Cells.Item(98, 9).Value = "unc-3 is a transcriptional factor"
$ws.Cells.Item(99, 9).Value = "unc-3 is a transcriptional factor"
$ws.Cells.Item(100, 9).Value = "unc-3 is a transcriptional factor"

# --- Remaining columns: A (date), B/C/D (ids), E (genotype "WT").
for ($i = 0; $i -lt 4; $i++) {
    $row = 97 + $i
    $ws.Cells.Item($row, 1).Value = 43322
    $ws.Cells.Item($row, 2).Value = 6625
    $ws.Cells.Item($row, 3).Value = 6641
    $ws.Cells.Item($row, 4).Value = 72762 + $i
    $ws.Cells.Item($row, 5).Value = "WT"
}

# --- Formatting: copy the existing centered-date style (column A, e.g. A6)
# and the existing centered general style (column B, e.g. B6) onto the new
# rows so they reuse the workbook's existing styles instead of creating new
# ones.
$ws.Range("A6").Copy()
$ws.Range("A97:A100").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("B97:I100").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Move the selection to where the author ended up after entering the
# data (one row below the last data row, column G).
$ws.Range("G117").Select() | Out-Null
